# Full implementation of passers. Not fully tested
#
# Adds a "Filas"/"Columnas" helper table (rows 9-17, columns R:Y and AA:AH)
# to the right of the existing chess-square data, and updates the active
# sheet's selection to the new Columnas block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 headers
$ws.Range("R9").Value = "Filas"
$ws.Range("AA9").Value = "Columnas"

# Rows 10-17: "Filas" block (columns R:Y, same value repeated across the row,
# counting down from 7 to 0) and "Columnas" block (columns AA:AH, 0..7 on
# every row).
for ($row = 10; $row -le 17; $row++) {
    $filas = 17 - $row

    for ($col = 18; $col -le 25; $col++) {
        $ws.Cells.Item($row, $col).Value = $filas
    }

    for ($i = 0; $i -le 7; $i++) {
        $col = 27 + $i
        $ws.Cells.Item($row, $col).Value = $i
    }
}

# Update the visible selection to the new Columnas block (also clears the
# old topLeftCell scroll-freeze that pointed at A6).
[void]$ws.Range("AA10:AH17").Select()
